# Auto-generated edit script: updates Leve profit/price calculation columns
# (H, I, J, K, L, M, N) across multiple sheets, per refreshed market data.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 51 (G51=5486)
$ws.Range("H51").Value = 1928.2
$ws.Range("I51").Value = 1766.6666
$ws.Range("J51").Value = 1997.4286
$ws.Range("K51").Value = 1766.6666
$ws.Range("L51").Value = 1997.4286
$ws.Range("M51").Value = -1282.6666
$ws.Range("N51").Value = -2965.4286

# Row 138 (G138=44169)
$ws.Range("H138").Value = 5925.205
$ws.Range("I138").Value = 3155.2942
$ws.Range("J138").Value = 6697.1475
$ws.Range("K138").Value = 9465.882599999999
$ws.Range("L138").Value = 20091.4425
$ws.Range("N138").Value = -30371.4425
$ws.Range("M138").Value = -4325.882599999999

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G32=44147)
$ws.Range("H32").Value = 33445.953
$ws.Range("I32").Value = 13376.16
$ws.Range("J32").Value = 117070.086
$ws.Range("K32").Value = 13376.16
$ws.Range("L32").Value = 117070.086
$ws.Range("M32").Value = -13089.16
$ws.Range("N32").Value = -117644.086

# Row 45 (G45=27714)
$ws.Range("H45").Value = 3688.32
$ws.Range("I45").Value = 2682
$ws.Range("J45").Value = 5826.75
$ws.Range("K45").Value = 2682
$ws.Range("L45").Value = 5826.75
$ws.Range("M45").Value = -2305
$ws.Range("N45").Value = -6580.75

# Row 55 (G55=2830)
$ws.Range("H55").Value = 35053
$ws.Range("J55").Value = 35053
$ws.Range("L55").Value = 35053
$ws.Range("N55").Value = -35683

# Row 80 (G80=10667)
$ws.Range("H80").Value = 37927
$ws.Range("J80").Value = 37927
$ws.Range("L80").Value = 37927
$ws.Range("N80").Value = -39923

# Row 83 (G83=10667)
$ws.Range("H83").Value = 37927
$ws.Range("J83").Value = 37927
$ws.Range("L83").Value = 113781
$ws.Range("N83").Value = -123765

# Row 103 (G103=18533)
$ws.Range("H103").Value = 41277.6
$ws.Range("J103").Value = 41277.6
$ws.Range("L103").Value = 41277.6
$ws.Range("N103").Value = -43621.6

# Row 119 (G119=26287)
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("N119").ClearContents()

# Row 131 (G131=34706)
$ws.Range("H131").Value = 57592.332
$ws.Range("J131").Value = 57592.332
$ws.Range("L131").Value = 57592.332
$ws.Range("N131").Value = -67672.33199999999

$ws = $wb.Worksheets.Item("BSM")
# Row 35 (G35=2350)
$ws.Range("H35").Value = 35474
$ws.Range("J35").Value = 35474
$ws.Range("L35").Value = 35474
$ws.Range("N35").Value = -36094

# Row 82 (G82=11877)
$ws.Range("H82").Value = 26212.5
$ws.Range("J82").Value = 42901
$ws.Range("L82").Value = 42901
$ws.Range("N82").Value = -43667

# Row 85 (G85=11877)
$ws.Range("H85").Value = 26212.5
$ws.Range("J85").Value = 42901
$ws.Range("L85").Value = 42901
$ws.Range("N85").Value = -45553

# Row 122 (G122=34096)
$ws.Range("H122").Value = 0
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

# Row 124 (G124=34245)
$ws.Range("H124").Value = 42375
$ws.Range("J124").Value = 42375
$ws.Range("L124").Value = 42375
$ws.Range("N124").Value = -52195

# Row 126 (G126=34398)
$ws.Range("H126").Value = 27815.455
$ws.Range("J126").Value = 27815.455
$ws.Range("L126").Value = 27815.455
$ws.Range("N126").Value = -37695.455

# Row 130 (G130=34682)
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

# Row 131 (G131=35396)
$ws.Range("H131").Value = 99980
$ws.Range("J131").Value = 99980
$ws.Range("L131").Value = 99980
$ws.Range("N131").Value = -110060

$ws = $wb.Worksheets.Item("CRP")
# Row 31 (G31=44023)
$ws.Range("H31").Value = 4189.609
$ws.Range("I31").Value = 2053.5173
$ws.Range("J31").Value = 7833.5293
$ws.Range("K31").Value = 2053.5173
$ws.Range("L31").Value = 7833.5293
$ws.Range("M31").Value = -1758.5173
$ws.Range("N31").Value = -8423.5293

# Row 34 (G34=44023)
$ws.Range("H34").Value = 4189.609
$ws.Range("I34").Value = 2053.5173
$ws.Range("J34").Value = 7833.5293
$ws.Range("K34").Value = 2053.5173
$ws.Range("L34").Value = 7833.5293
$ws.Range("M34").Value = -1851.5173
$ws.Range("N34").Value = -8237.5293

# Row 41 (G41=1917)
$ws.Range("H41").Value = 15997.6
$ws.Range("J41").Value = 29965
$ws.Range("L41").Value = 29965
$ws.Range("N41").Value = -30821

# Row 50 (G50=1862)
$ws.Range("H50").Value = 7864.875
$ws.Range("J50").Value = 8976.571
$ws.Range("L50").Value = 8976.571
$ws.Range("N50").Value = -10226.571

# Row 51 (G51=2039)
$ws.Range("H51").Value = 9386.833000000001
$ws.Range("J51").Value = 9386.833000000001
$ws.Range("L51").Value = 9386.833000000001
$ws.Range("N51").Value = -10858.833

# Row 60 (G60=1937)
$ws.Range("H60").Value = 12086.75
$ws.Range("J60").Value = 12086.75
$ws.Range("L60").Value = 12086.75
$ws.Range("N60").Value = -13108.75

# Row 61 (G61=2039)
$ws.Range("H61").Value = 9386.833000000001
$ws.Range("J61").Value = 9386.833000000001
$ws.Range("L61").Value = 9386.833000000001
$ws.Range("N61").Value = -10082.833

$ws = $wb.Worksheets.Item("CUL")
# Row 92 (G92=19841)
$ws.Range("H92").Value = 2500
$ws.Range("J92").Value = 2500
$ws.Range("L92").Value = 7500
$ws.Range("N92").Value = -9996

# Row 113 (G113=27843)
$ws.Range("H113").Value = 1058.6666
$ws.Range("I113").Value = 676.6667
$ws.Range("J113").Value = 1186
$ws.Range("K113").Value = 2030.0001
$ws.Range("L113").Value = 3558
$ws.Range("M113").Value = 139.9999
$ws.Range("N113").Value = -7898

# Row 131 (G131=36060)
$ws.Range("H131").Value = 888.9
$ws.Range("J131").Value = 888.9
$ws.Range("L131").Value = 2666.7
$ws.Range("N131").Value = -12746.7

# Row 137 (G137=44088)
$ws.Range("H137").Value = 2377.389
$ws.Range("I137").Value = 2104.6155
$ws.Range("J137").Value = 3086.6
$ws.Range("K137").Value = 6313.8465
$ws.Range("L137").Value = 9259.799999999999
$ws.Range("M137").Value = -1213.8465
$ws.Range("N137").Value = -19459.8

$ws = $wb.Worksheets.Item("GSM")
# Row 57 (G57=2876)
$ws.Range("H57").Value = 16365.1875
$ws.Range("J57").Value = 15896.2
$ws.Range("L57").Value = 15896.2
$ws.Range("N57").Value = -17536.2

# Row 117 (G117=26185)
$ws.Range("H117").Value = 14600
$ws.Range("J117").Value = 14600
$ws.Range("L117").Value = 14600
$ws.Range("N117").Value = -21484

$ws = $wb.Worksheets.Item("LTW")
# Row 95 (G95=18221)
$ws.Range("H95").Value = 34672
$ws.Range("J95").Value = 34672
$ws.Range("L95").Value = 34672
$ws.Range("N95").Value = -40164

# Row 136 (G136=44060)
$ws.Range("H136").Value = 1708.2
$ws.Range("I136").Value = 1509.1111
$ws.Range("J136").Value = 3500
$ws.Range("K136").Value = 4527.3333
$ws.Range("L136").Value = 10500
$ws.Range("M136").Value = -1977.3333
$ws.Range("N136").Value = -15600

$ws = $wb.Worksheets.Item("WVR")
# Row 98 (G98=18374)
$ws.Range("H98").Value = 43918.75
$ws.Range("J98").Value = 43918.75
$ws.Range("L98").Value = 43918.75
$ws.Range("N98").Value = -49908.75

# Row 109 (G109=27161)
$ws.Range("H109").Value = 21033.666
$ws.Range("J109").Value = 21033.666
$ws.Range("L109").Value = 21033.666
$ws.Range("N109").Value = -23807.666
